$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.29903
$ws.Range("H2").Value = 6.89709
$ws.Range("I2").Value = 0.04075801785348079
$ws.Range("J2").Value = 0.04075801785348079
$ws.Range("M2").Value = 14.11187666666667
$ws.Range("N2").Value = 42.33562999999999
$ws.Range("O2").Value = 0.08862966207485527
$ws.Range("P2").Value = 0.08862966207485526
$ws.Range("Q2").Value = 32.44362781296667
$ws.Range("R2").Value = 291.9926503167
$ws.Range("S2").Value = 0.00361236934919492
$ws.Range("T2").Value = 0.00361236934919492
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.29903
$ws.Range("H3").Value = 6.89709
$ws.Range("I3").Value = 0.04075801785348079
$ws.Range("J3").Value = 0.04075801785348079
$ws.Range("O3").Value = 0.7176943460983047
$ws.Range("P3").Value = 0.7176943460983046
$ws.Range("Q3").Value = 262.7180077547633
$ws.Range("R3").Value = 2364.46206979287
$ws.Range("S3").Value = 0.02925179897161692
$ws.Range("T3").Value = 0.02925179897161692
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.29903
$ws.Range("H4").Value = 6.89709
$ws.Range("I4").Value = 0.04075801785348079
$ws.Range("J4").Value = 0.04075801785348079
$ws.Range("O4").Value = 0.1936759918268401
$ws.Range("P4").Value = 0.1936759918268401
$ws.Range("Q4").Value = 70.89671389957667
$ws.Range("R4").Value = 638.07042509619
$ws.Range("S4").Value = 0.00789384953266895
$ws.Range("T4").Value = 0.00789384953266895
# Row 5
$ws.Range("I5").Value = 0.3949230674234065
$ws.Range("J5").Value = 0.3949230674234066
$ws.Range("M5").Value = 14.11187666666667
$ws.Range("N5").Value = 42.33562999999999
$ws.Range("O5").Value = 0.08862966207485527
$ws.Range("P5").Value = 0.08862966207485526
$ws.Range("Q5").Value = 314.3611414151711
$ws.Range("R5").Value = 2829.25027273654
$ws.Range("S5").Value = 0.03500189801130181
$ws.Range("T5").Value = 0.0350018980113018
# Row 6
$ws.Range("I6").Value = 0.3949230674234065
$ws.Range("J6").Value = 0.3949230674234066
$ws.Range("O6").Value = 0.7176943460983047
$ws.Range("P6").Value = 0.7176943460983046
$ws.Range("S6").Value = 0.2834340526335785
$ws.Range("T6").Value = 0.2834340526335785
# Row 7
$ws.Range("I7").Value = 0.3949230674234065
$ws.Range("J7").Value = 0.3949230674234066
$ws.Range("O7").Value = 0.1936759918268401
$ws.Range("P7").Value = 0.1936759918268401
$ws.Range("S7").Value = 0.07648711677852632
$ws.Range("T7").Value = 0.07648711677852632
# Row 8
$ws.Range("I8").Value = 0.5643189147231126
$ws.Range("J8").Value = 0.5643189147231126
$ws.Range("M8").Value = 14.11187666666667
$ws.Range("N8").Value = 42.33562999999999
$ws.Range("O8").Value = 0.08862966207485527
$ws.Range("P8").Value = 0.08862966207485526
$ws.Range("Q8").Value = 449.2012566192633
$ws.Range("R8").Value = 4042.811309573369
$ws.Range("S8").Value = 0.05001539471435854
$ws.Range("T8").Value = 0.05001539471435853
# Row 9
$ws.Range("I9").Value = 0.5643189147231126
$ws.Range("J9").Value = 0.5643189147231126
$ws.Range("O9").Value = 0.7176943460983047
$ws.Range("P9").Value = 0.7176943460983046
$ws.Range("S9").Value = 0.4050084944931093
$ws.Range("T9").Value = 0.4050084944931092
# Row 10
$ws.Range("I10").Value = 0.5643189147231126
$ws.Range("J10").Value = 0.5643189147231126
$ws.Range("O10").Value = 0.1936759918268401
$ws.Range("P10").Value = 0.1936759918268401
$ws.Range("S10").Value = 0.1092950255156449
$ws.Range("T10").Value = 0.1092950255156448

Write-Host "Applied 82 cell updates"
